$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = -7.560999999999998
$ws.Range("D18").Value = -8.543000000000003
$ws.Range("D20").Value = -7.62
$ws.Range("D27").Value = -8.513000000000002
$ws.Range("D69").Value = -7.204999999999998
$ws.Range("D76").Value = -7.753000000000002
$ws.Range("D82").Value = -8.399000000000001
